$wb = $excel.ActiveWorkbook

# Week 16 logged stats - update row 2 ("H" row) on both OFF and DEF sheets

$wsOFF = $wb.Worksheets.Item("OFF")
$wsOFF.Range("B2").Value = 439
$wsOFF.Range("C2").Value = 325
$wsOFF.Range("D2").Value = 110
$wsOFF.Range("E2").Value = 50

$wsDEF = $wb.Worksheets.Item("DEF")
$wsDEF.Range("B2").Value = 434
$wsDEF.Range("C2").Value = 302
$wsDEF.Range("D2").Value = 85
$wsDEF.Range("E2").Value = 35
